# Auto-generated edit script: applies market/profit data refresh
# to the 8 crafting-class worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 2225.4
$ws.Range("I11").Value = 2225.4
$ws.Range("K11").Value = 2225.4
$ws.Range("M11").Value = -2085.4

$ws.Range("H32").Value = 4173.5264
$ws.Range("I32").Value = 4037.375
$ws.Range("J32").Value = 4272.5454
$ws.Range("K32").Value = 4037.375
$ws.Range("L32").Value = 4272.5454
$ws.Range("M32").Value = -3711.375
$ws.Range("N32").Value = -4924.5454

$ws.Range("H55").Value = 606
$ws.Range("J55").Value = 1
$ws.Range("L55").Value = 1
$ws.Range("N55").Value = -429

$ws.Range("H62").Value = 16155.333
$ws.Range("J62").Value = 17612.375
$ws.Range("L62").Value = 17612.375
$ws.Range("N62").Value = -18860.375

$ws.Range("H65").Value = 16155.333
$ws.Range("J65").Value = 17612.375
$ws.Range("L65").Value = 88061.875
$ws.Range("N65").Value = -94301.875

$ws.Range("H96").Value = 2351.4614
$ws.Range("J96").Value = 2824.2
$ws.Range("L96").Value = 8472.599999999999
$ws.Range("N96").Value = -11218.6

$ws.Range("H100").Value = 5132.385
$ws.Range("I100").Value = 1358.5454
$ws.Range("K100").Value = 1358.5454
$ws.Range("M100").Value = -817.5454

$ws.Range("H131").Value = 3539.389
$ws.Range("I131").Value = 1973.5454
$ws.Range("K131").Value = 5920.6362
$ws.Range("M131").Value = -880.6361999999999

$ws.Range("H137").Value = 7404.6816
$ws.Range("I137").Value = 6745.25
$ws.Range("J137").Value = 13999
$ws.Range("K137").Value = 20235.75
$ws.Range("L137").Value = 41997
$ws.Range("M137").Value = -17685.75
$ws.Range("N137").Value = -47097

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1013.0476
$ws.Range("I2").Value = 946.2353000000001
$ws.Range("K2").Value = 946.2353000000001
$ws.Range("M2").Value = -833.2353000000001

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H61").Value = 41674704
$ws.Range("I61").Value = 45461016
$ws.Range("J61").Value = 35724784
$ws.Range("K61").Value = 45461016
$ws.Range("L61").Value = 35724784
$ws.Range("M61").Value = -45460804
$ws.Range("N61").Value = -35725208

$ws.Range("H92").Value = 55731.332
$ws.Range("J92").Value = 55731.332
$ws.Range("L92").Value = 55731.332
$ws.Range("N92").Value = -60723.332

$ws.Range("H102").Value = 18483.555
$ws.Range("I102").Value = 18483.555
$ws.Range("K102").Value = 18483.555
$ws.Range("M102").Value = -16861.555

$ws.Range("H110").Value = 1114.95
$ws.Range("I110").Value = 927
$ws.Range("K110").Value = 927
$ws.Range("M110").Value = 1118

$ws.Range("H116").Value = 1013.0476
$ws.Range("I116").Value = 946.2353000000001
$ws.Range("K116").Value = 946.2353000000001
$ws.Range("M116").Value = 1347.7647

$ws.Range("H122").Value = 3292.25
$ws.Range("I122").Value = 1941.1818
$ws.Range("K122").Value = 5823.5454
$ws.Range("M122").Value = -3373.5454

$ws.Range("H136").Value = 41674704
$ws.Range("I136").Value = 45461016
$ws.Range("J136").Value = 35724784
$ws.Range("K136").Value = 136383048
$ws.Range("L136").Value = 107174352
$ws.Range("M136").Value = -136380498
$ws.Range("N136").Value = -107179452

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1013.0476
$ws.Range("I3").Value = 946.2353000000001
$ws.Range("K3").Value = 946.2353000000001
$ws.Range("M3").Value = -832.2353000000001

$ws.Range("H22").Value = 913.3333
$ws.Range("I22").Value = 652.25
$ws.Range("K22").Value = 652.25
$ws.Range("M22").Value = -479.25

$ws.Range("H99").Value = 2233.8838
$ws.Range("I99").Value = 1781.3939
$ws.Range("K99").Value = 1781.3939
$ws.Range("M99").Value = -283.3939

$ws.Range("H134").Value = 918846.3
$ws.Range("I134").Value = 1914.2
$ws.Range("K134").Value = 5742.6
$ws.Range("M134").Value = -3207.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 703.4286
$ws.Range("J16").Value = 1606
$ws.Range("L16").Value = 1606
$ws.Range("N16").Value = -2180

$ws.Range("H58").Value = 5635.5557
$ws.Range("I58").Value = 4186.8335
$ws.Range("K58").Value = 4186.8335
$ws.Range("M58").Value = -3983.8335

$ws.Range("H92").Value = 69530
$ws.Range("J92").Value = 69530
$ws.Range("L92").Value = 69530
$ws.Range("N92").Value = -74522

$ws.Range("H113").Value = 703.4286
$ws.Range("J113").Value = 1606
$ws.Range("L113").Value = 1606
$ws.Range("N113").Value = -5946

$ws.Range("H122").Value = 1642.5
$ws.Range("I122").Value = 1390.125
$ws.Range("K122").Value = 4170.375
$ws.Range("M122").Value = -1720.375

$ws.Range("H136").Value = 5635.5557
$ws.Range("I136").Value = 4186.8335
$ws.Range("K136").Value = 12560.5005
$ws.Range("M136").Value = -10010.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 118.84091
$ws.Range("I2").Value = 73.25
$ws.Range("J2").Value = 135.9375
$ws.Range("K2").Value = 439.5
$ws.Range("L2").Value = 815.625
$ws.Range("M2").Value = -326.5
$ws.Range("N2").Value = -1041.625

$ws.Range("H5").Value = 2389.3635
$ws.Range("I5").Value = 2328.3
$ws.Range("K5").Value = 6984.900000000001
$ws.Range("M5").Value = -6872.900000000001

$ws.Range("H60").Value = 2637.158
$ws.Range("I60").Value = 1175
$ws.Range("J60").Value = 2809.1765
$ws.Range("K60").Value = 3525
$ws.Range("L60").Value = 8427.529500000001
$ws.Range("M60").Value = -3274
$ws.Range("N60").Value = -8929.529500000001

$ws.Range("H104").Value = 4000
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()

$ws.Range("H121").Value = 1213.091
$ws.Range("I121").Value = 514.6923
$ws.Range("J121").Value = 2221.889
$ws.Range("K121").Value = 1544.0769
$ws.Range("L121").Value = 6665.667
$ws.Range("M121").Value = -234.0769
$ws.Range("N121").Value = -9285.667000000001

$ws.Range("H131").Value = 254861.62
$ws.Range("I131").Value = 2501197.5
$ws.Range("J131").Value = 18405.21
$ws.Range("K131").Value = 7503592.5
$ws.Range("L131").Value = 55215.63
$ws.Range("M131").Value = -7498552.5
$ws.Range("N131").Value = -65295.63

$ws.Range("H135").Value = 2389.3635
$ws.Range("I135").Value = 2328.3
$ws.Range("K135").Value = 20954.7
$ws.Range("M135").Value = -18419.7

$ws.Range("H137").Value = 4355.5386
$ws.Range("I137").Value = 4654.6665
$ws.Range("J137").Value = 4099.143
$ws.Range("K137").Value = 13963.9995
$ws.Range("L137").Value = 12297.429
$ws.Range("M137").Value = -8863.999500000002
$ws.Range("N137").Value = -22497.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 1750
$ws.Range("I28").Value = 2500
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 2500
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = -2308
$ws.Range("N28").Value = -1384

$ws.Range("H92").Value = 41750.332
$ws.Range("J92").Value = 41750.332
$ws.Range("L92").Value = 41750.332
$ws.Range("N92").Value = -45494.332

$ws.Range("H102").Value = 3121.0908
$ws.Range("I102").Value = 2477
$ws.Range("K102").Value = 2477
$ws.Range("M102").Value = -855

$ws.Range("H132").Value = 142874270
$ws.Range("I132").Value = 200001980
$ws.Range("K132").Value = 600005940
$ws.Range("M132").Value = -600003410

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4766.381
$ws.Range("I40").Value = 4284
$ws.Range("K40").Value = 4284
$ws.Range("M40").Value = -4148

$ws.Range("H55").Value = 38462052
$ws.Range("I55").Value = 62500524
$ws.Range("K55").Value = 62500524
$ws.Range("M55").Value = -62500351

$ws.Range("H132").Value = 1836916.5
$ws.Range("I132").Value = 252873.5
$ws.Range("J132").Value = 5005002.5
$ws.Range("K132").Value = 758620.5
$ws.Range("L132").Value = 15015007.5
$ws.Range("M132").Value = -756090.5
$ws.Range("N132").Value = -15020067.5

$ws.Range("H136").Value = 127961.84
$ws.Range("I136").Value = 26099.8
$ws.Range("K136").Value = 78299.39999999999
$ws.Range("M136").Value = -75749.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H96").Value = 1768.6
$ws.Range("J96").Value = 1768.6
$ws.Range("L96").Value = 1768.6
$ws.Range("N96").Value = -4514.6

$ws.Range("H126").Value = 2690.6
$ws.Range("I126").Value = 2690.6
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8071.799999999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5601.799999999999
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 1005340.1
$ws.Range("I132").Value = 6149.3335
$ws.Range("K132").Value = 18448.0005
$ws.Range("M132").Value = -15918.0005

$ws.Range("H136").Value = 6669.2666
$ws.Range("I136").Value = 7436.5
$ws.Range("K136").Value = 22309.5
$ws.Range("M136").Value = -19759.5
